# The deck originally carried the "Integral" theme on its (visible) slide
# master and the default Office "Office Theme" colors on its notes master.
# The authored change swaps the two: the slide master's theme becomes the
# standard Office palette (what used to live behind the notes master),
# while the notes master keeps the palette that used to drive the slides.
#
# PowerPoint's object model only exposes the *active* (slide-master-facing)
# theme's ThemeColorScheme for writing, so we drive the swap by pushing the
# target "Office Theme" RGB values onto that scheme, one ThemeColorScheme
# slot at a time, in the standard dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink
# order.

function HexToBgr($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $b * 65536 + $g * 256 + $r
}

$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$scheme = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Length; $i++) {
    $scheme.Item($i).RGB = HexToBgr($officeThemeColors[$i - 1])
}
